$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 6 ---
$ws.Range("A6").Value = "Sim getauscht"
$ws.Range("B6").Value = 2
$ws.Range("I6").Value = "'+COPS: 0,0,`"Telekom.de 1nce.net`",9"
$ws.Range("G6").Value = "Geht online nach ca 40 runden"
$linkC6 = "https://github.com/Qrist0ph/pirmicboard_david/blob/a2ef3bf7a16e352b20008353330c1dd8c7383065/unittests/7080gconnect/main.py"
$ws.Hyperlinks.Add($ws.Range("C6"), $linkC6, "", "", $linkC6) | Out-Null
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("C6").Value = "pirmicboard_david/unittests/7080gconnect/main.py at a2ef3bf7a16e352b20008353330c1dd8c7383065 · Qrist0ph/pirmicboard_david (github.com)"
$ws.Range("D6").Value = "'8988228066603839868"

# --- Row 7 ---
$ws.Range("A7").Value = "Sim getauscht"
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = "'8988228066603839867"
$ws.Range("G7").Value = "Online nach runde 27"

# --- Row 8 ---
$ws.Range("A8").Value = "Board getauscht"
$ws.Range("B8").Value = 1
$ws.Range("D8").Value = "'8988228066603839868"

# --- Hyperlinks for C7 / C8 (created in this order so the shared string
# for the 5e8b1c0... link text is interned before the G8 text) ---
$linkC7C8 = "https://github.com/Qrist0ph/pirmicboard_david/blob/5e8b1c0b49c3b408311a3a79990f6b82077685f9/unittests/7080gconnect/main.py"
$descC7C8 = "pirmicboard_david/unittests/7080gconnect/main.py at 5e8b1c0b49c3b408311a3a79990f6b82077685f9 · Qrist0ph/pirmicboard_david (github.com)"
$ws.Hyperlinks.Add($ws.Range("C7"), $linkC7C8, "", "", $linkC7C8) | Out-Null
$ws.Range("C7").Style = "Hyperlink"
$ws.Range("C7").Value = $descC7C8
$ws.Hyperlinks.Add($ws.Range("C8"), $linkC7C8, "", "", $linkC7C8) | Out-Null
$ws.Range("C8").Style = "Hyperlink"
$ws.Range("C8").Value = $descC7C8

# --- Row 8 continued ---
$ws.Range("G8").Value = "ein reset war notwendig, dann online gegangen"

# --- Selection ---
$ws.Range("G9").Select() | Out-Null
